$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Supplier fix review comments:
# Row 30 (CRS_supplier_004): drop the "_[Color]_" segment from the product id convention
$ws.Range("B30").Value = "Product Id follow ID convention [Brand/Category]_[ProductName]_[Size/Version]"

# Row 27 (CRS_supplier_001): clarify supported platforms for "add product"
$ws.Range("B27").Value = "Supplier can add product with data (product id,product photo, product price, product version, product platform [IOS&Android] )"

# Row 31 (CRS_supplier_005): product photo & platform become mandatory, version rewording
$ws.Range("B31").Value = "If the supplier wants to add a new product the product data which is mandatory (product id, product price , product version, product photo& product platform)"

# Move the active selection/scroll position to reflect the reviewed rows
$ws.Range("B33").Select()
